# Apply updates to column F (dSF) values for specific rows on Sheet1
# as part of a data repull / mean calculation refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -6
    9  = -4
    11 = -1
    12 = -3
    13 = -3
    14 = -3
    15 = -1
    16 = 1
    18 = 5
    19 = -3
    20 = -1
    21 = -1
    22 = 1
    23 = 4
    24 = -2
    25 = 1
    26 = -1
    28 = 4
    30 = 3
    31 = 4
    34 = 0
    35 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
